$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix spreadsheet headers (GPLIM-2588)
$ws.Range("A1").Value = "Specimen_Number"
$ws.Range("F1").Value = "SAMPLE_TYPE"

# Update the active selection to match the saved view state
$ws.Range("F1").Select()
